$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.153.38'
$ws.Range("E2").Value = '  -0.72%  '

$ws.Range("D3").Value = '2.468.08'
$ws.Range("E3").Value = '  -2.22%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").Value = '''552.14'
$ws.Range("E5").Value = '  -1.64%  '

$ws.Range("D6").Value = '''147.31'
$ws.Range("E6").Value = '  -1.17%  '

$ws.Range("D8").Value = '''0.589'
$ws.Range("E8").Value = '  -3.82%  '

$ws.Range("D9").Value = '2.467.17'
$ws.Range("E9").Value = '  -2.25%  '

$ws.Range("D10").Value = '''0.107'
$ws.Range("E10").Value = '  -3.41%  '

$ws.Range("E11").Value = '  -0.06%  '

$ws.Range("D12").Value = '''5.43'
$ws.Range("E12").Value = '  -1.22%  '

$ws.Range("D13").Value = '''0.353'
$ws.Range("E13").Value = '  -3.65%  '

$ws.Range("D14").Value = '''26.31'
$ws.Range("E14").Value = '  -1.73%  '

$ws.Range("D15").Value = '2.918.30'
$ws.Range("E15").Value = '  -2.07%  '

$ws.Range("D16").Value = '''0.0000169'
$ws.Range("E16").Value = '  -0.46%  '

$ws.Range("D17").Value = '62.124.51'
$ws.Range("E17").Value = '  -0.55%  '

$ws.Range("D18").Value = '2.472.13'
$ws.Range("E18").Value = '  -1.62%  '

$ws.Range("D19").Value = '''10.98'
$ws.Range("E19").Value = '  -3.88%  '

$ws.Range("D20").Value = '''7.07'
$ws.Range("E20").Value = '  -0.65%  '

$ws.Range("D21").Value = '''4.18'
$ws.Range("E21").Value = '  -2.63%  '

$ws.Range("D22").Value = '''321.85'
$ws.Range("E22").Value = '  -1.91%  '

$ws.Range("E23").Value = '  +0.14%  '

$ws.Range("D24").Value = '''1.88'
$ws.Range("E24").Value = '  +4.81%  '

$ws.Range("D25").Value = '''64.09'
$ws.Range("E25").Value = '  -1.96%  '

$ws.Range("D26").Value = '0.0₃0993'
$ws.Range("E26").Value = '  -7.80%  '

$ws.Range("D27").Value = '2.590.05'
$ws.Range("E27").Value = '  -2.27%  '

$ws.Range("B28").Value = 'Fetch.AI'
$ws.Range("C28").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D28").Value = '''1.49'
$ws.Range("E28").Value = '  -4.88%  '

$ws.Range("B29").Value = 'Binance-PegBSC-USD'
$ws.Range("C29").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D29").Value = '''0.999'
$ws.Range("E29").Value = '  -0.06%  '

$ws.Range("D30").Value = '''536.74'
$ws.Range("E30").Value = '  -3.17%  '

$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").Value = '''8.30'
$ws.Range("E31").Value = '  -4.94%  '

$ws.Range("B32").Value = 'Aptos'
$ws.Range("C32").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D32").Value = '''7.79'
$ws.Range("E32").Value = '  -2.64%  '

$ws.Range("D33").Value = '''0.148'
$ws.Range("E33").Value = '  -4.92%  '

$ws.Range("D34").Value = '''1.90'
$ws.Range("E34").Value = '  -2.64%  '

$ws.Range("D35").Value = '''1.63'
$ws.Range("E35").Value = '  +0.30%  '

$ws.Range("D36").Value = '''5.75'
$ws.Range("E36").Value = '  -6.57%  '

$ws.Range("E37").Value = '  +0.06%  '

$ws.Range("D38").Value = '''4.83'
$ws.Range("E38").Value = '  -3.13%  '

$ws.Range("D39").Value = '''0.382'
$ws.Range("E39").Value = '  -1.03%  '

$ws.Range("D40").Value = '''18.31'
$ws.Range("E40").Value = '  -3.27%  '

$ws.Range("D41").Value = '''1.77'
$ws.Range("E41").Value = '  +1.62%  '

$ws.Range("D42").Value = '''139.90'
$ws.Range("E42").Value = '  -6.66%  '

$ws.Range("E43").Value = '  +0.14%  '

$ws.Range("D44").Value = '''40.47'
$ws.Range("E44").Value = '  -1.32%  '

$ws.Range("D45").Value = '''2.32'
$ws.Range("E45").Value = '  -3.94%  '

$ws.Range("D46").Value = '''144.86'
$ws.Range("E46").Value = '  -4.06%  '

$ws.Range("D47").Value = '''3.63'
$ws.Range("E47").Value = '  -2.06%  '

$ws.Range("D48").Value = '''21.80'
$ws.Range("E48").Value = '  -1.24%  '

$ws.Range("D49").Value = '''0.0534'
$ws.Range("E49").Value = '  -3.22%  '

$ws.Range("D50").Value = '''0.593'
$ws.Range("E50").Value = '  -1.72%  '

$ws.Range("D51").Value = '''0.0935'
$ws.Range("E51").Value = '  -2.94%  '
